# "fjernet mellomrom" - removed the trailing space from the business name
# "Søya hytteutleie " (shared string) in cell E4 of the
# overnatting-i-gjesdal-kommune sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "Søya hytteutleie"

# Leave the selection where the author's last save left it.
[void]$ws.Range("E5").Select()
